# Regional Availability Factor.xlsx -- "updated 4.0 files and mdl"
#
# Substantive changes in this revision:
#  1. About sheet: bump the "last updated" date stamp (2024-03-15 -> 2024-03-28).
#  2. RAF-capacity sheet: raise the capacity-credit multiplier for the two
#     hydrogen technologies (hydrogen combustion turbine / hydrogen combined
#     cycle) from 0.3 to 1.
#  3. View state: RAF-capacity becomes the active/selected sheet (instead of
#     RAF-generation), scrolled/zoomed to show the bottom of the list with the
#     hydrogen rows selected, and its first column is narrowed to fit the
#     longer row labels.

$wb = $excel.ActiveWorkbook

# --- 1. About sheet: refresh the date stamp in C1 -------------------------
$about = $wb.Worksheets.Item("About")
$about.Range("C1").Value = 45379   # serial date for 2024-03-28

# --- 2. RAF-capacity sheet: hydrogen RAF values 0.3 -> 1 -------------------
$cap = $wb.Worksheets.Item("RAF-capacity")
$cap.Range("B24").Value = 1   # hydrogen combustion turbine
$cap.Range("B25").Value = 1   # hydrogen combined cycle

# Narrow column A slightly to better fit the labels (was default width).
$cap.Columns("A").ColumnWidth = 28.1666666666667

# --- 3. View state: make RAF-capacity the active sheet, scrolled/zoomed ---
$cap.Activate()
$cap.Range("B25").Select()

$win = $excel.ActiveWindow
$win.Zoom = 80
$win.ScrollRow = 14
$win.ScrollColumn = 1
